$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F values updated
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 385
$ws1.Range("F3").Value = 1010
$ws1.Range("F4").Value = 242
$ws1.Range("F5").Value = 1368
$ws1.Range("F6").Value = 8428
$ws1.Range("F8").Value = 474
$ws1.Range("F9").Value = 621
$ws1.Range("F10").Value = 236
$ws1.Range("F11").Value = 143
$ws1.Range("F12").Value = 3383
$ws1.Range("F15").Value = 51
$ws1.Range("F16").Value = 912
$ws1.Range("F17").Value = 139
$ws1.Range("F21").Value = 1988

# Sheet "全部类型" (all types) - column F values updated
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 386
$ws4.Range("F3").Value = 1010
$ws4.Range("F4").Value = 242
$ws4.Range("F5").Value = 1368
$ws4.Range("F6").Value = 8428
$ws4.Range("F8").Value = 474
$ws4.Range("F9").Value = 621
$ws4.Range("F10").Value = 236
$ws4.Range("F11").Value = 143
$ws4.Range("F12").Value = 3383
$ws4.Range("F15").Value = 51
$ws4.Range("F16").Value = 912
$ws4.Range("F17").Value = 139
$ws4.Range("F21").Value = 1988
